# [ADDITIONAL SCRAPING] added scraping code for extra bowling attributes and excel sheets
#
# 1) Clean up "ODI Batting Extra": drop cells that only ever held an empty
#    placeholder string (no real scraped value) so they are no longer present.
# 2) Add a brand-new "ODI Bowling Extra" sheet (the bowling counterpart of the
#    existing "ODI Batting Extra" sheet) with MATCH_CODE / MAIDEN_OVERS /
#    PERCENT_WICKETS_OF_ALL columns and 20 rows of scraped data.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Step 1: "ODI Batting Extra" - remove leftover empty placeholder cells
# ---------------------------------------------------------------------------
$battingExtra = $wb.Worksheets.Item("ODI Batting Extra")

$emptyCells = @(
    "C2", "D2", "E2",
    "C3", "D3", "E3",
    "E6",
    "E7",
    "B9", "C9", "D9", "E9",
    "B14", "C14", "D14", "E14",
    "B15", "C15", "D15", "E15",
    "C16", "D16", "E16",
    "E18",
    "E19",
    "C20", "D20", "E20"
)

foreach ($addr in $emptyCells) {
    $battingExtra.Range($addr).ClearContents()
}

# ---------------------------------------------------------------------------
# Step 2: Add the new "ODI Bowling Extra" sheet at the end of the workbook
# ---------------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$bowlingExtra = $wb.Worksheets.Add($null, $lastSheet)
$bowlingExtra.Name = "ODI Bowling Extra"

# Header row (bold, thin border, centered - matches the other sheets' headers)
$header = $bowlingExtra.Range("A1:C1")
$header.Font.Bold = $true
$header.HorizontalAlignment = -4108
$header.VerticalAlignment = -4160
$header.Borders.LineStyle = 1

$bowlingExtra.Range("A1").Value = "MATCH_CODE"
$bowlingExtra.Range("B1").Value = "MAIDEN_OVERS"
$bowlingExtra.Range("C1").Value = "PERCENT_WICKETS_OF_ALL"

# Keep MATCH_CODE / MAIDEN_OVERS / percentage values as plain text, exactly as
# scraped (so "4113" / "0" / "10.00%" do not get auto-coerced into numbers).
$dataRange = $bowlingExtra.Range("A2:C21")
$dataRange.NumberFormat = "@"

$rows = @(
    @("4113", "1",  "20.00%"),
    @("4152", "0",  "20.00%"),
    @("4156", "0",  "10.00%"),
    @("4158", "0",  "20.00%"),
    @("4163", "0",  "10.00%"),
    @("4257", "1",  "10.00%"),
    @("4285", "0",  "20.00%"),
    @("4295", $null, $null),
    @("4301", "0",  $null),
    @("4391", "1",  "10.00%"),
    @("4394", "1",  "20.00%"),
    @("4397", "0",  "10.00%"),
    @("4426", $null, $null),
    @("4439", $null, $null),
    @("4442", "0",  $null),
    @("4444", "0",  "10.00%"),
    @("4446", "0",  "10.00%"),
    @("4466", "0",  $null),
    @("4467", "1",  "10.00%"),
    @("4468", "0",  "10.00%")
)

$r = 2
foreach ($row in $rows) {
    $bowlingExtra.Cells.Item($r, 1).Value = $row[0]
    if ($null -ne $row[1]) {
        $bowlingExtra.Cells.Item($r, 2).Value = $row[1]
    }
    if ($null -ne $row[2]) {
        $bowlingExtra.Cells.Item($r, 3).Value = $row[2]
    }
    $r = $r + 1
}
